$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 96, shifting existing rows 96-103 down to 97-104.
$ws.Rows.Item(96).Insert()

# Populate the new row 96 with the new weekly data record.
$ws.Cells.Item(96, 1).Value = 10
$ws.Cells.Item(96, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(96, 3).Value = "La Araucanía"
$ws.Cells.Item(96, 4).Value = 45106
$ws.Cells.Item(96, 5).Value = 9
$ws.Cells.Item(96, 6).Value = "Fruta"
$ws.Cells.Item(96, 7).Value = 100108
$ws.Cells.Item(96, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(96, 9).Value = 100108003
$ws.Cells.Item(96, 10).Value = "Maracuyá"
$ws.Cells.Item(96, 11).Value = "Sin especificar"
$ws.Cells.Item(96, 12).Value = "Primera"
$ws.Cells.Item(96, 13).Value = 40
$ws.Cells.Item(96, 14).Value = 50000
$ws.Cells.Item(96, 15).Value = 50000
$ws.Cells.Item(96, 16).Value = 50000
$ws.Cells.Item(96, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(96, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(96, 19).Value = 2778
$ws.Cells.Item(96, 20).Value = 18
